$d = $word.ActiveDocument

# The "Bibliografia" paragraph currently holds a single run whose w:t text
# concatenates four numbered references with no separation. The edit splits
# that single <w:t> into four <w:t> runs separated by <w:br/> elements
# (manual line breaks), one per numbered reference, while keeping everything
# inside the same <w:r>.
#
# Word's Find/Replace "^l" replacement code inserts a manual line break
# (<w:br/>) and naturally splits the surrounding text into separate <w:t>
# nodes within the same run, which matches the desired OOXML shape exactly.

$d.Content.Find.Execute(
    "2005.2. MONTGOMERY", $true, $false, $false, $false, $false,
    $true, 1, $false, "2005.^l2. MONTGOMERY", 2) | Out-Null

$d.Content.Find.Execute(
    "2004.3. GRANT", $true, $false, $false, $false, $false,
    $true, 1, $false, "2004.^l3. GRANT", 2) | Out-Null

$d.Content.Find.Execute(
    "1996.4. WERKENA", $true, $false, $false, $false, $false,
    $true, 1, $false, "1996.^l4. WERKENA", 2) | Out-Null
